# ---------------------------------------------------------------------------
# Applies the CHEESE invoice-proforma edit described by the commit:
#  - resize the two-column table (left col wider, right col narrower)
#  - tidy the "Payment" header spacing
#  - simplify the "(... see note below ...)" parenthetical
#  - broaden the householder/property-owner line
#  - rewrite the bold "Note:" paragraph with new payment rules + expenses blurb
#  - add a new hyperlink to the CHEESE energy-tracer-resources page
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Useful special characters (avoids any source-encoding ambiguity).
$pound  = [char]0x00A3   # £
$rsquo  = [char]0x2019   # '
$ndash  = [char]0x2013   # -

# ---------------------------------------------------------------------------
# 1) Table column widths: 5813/3685 (dxa) -> 6238/3260 (dxa); dxa = 20 * pt
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 6238 / 20.0
$t.Columns.Item(2).Width = 3260 / 20.0

# ---------------------------------------------------------------------------
# 2) "Payment" header cell: trim the leading run of spaces
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "                   Payment ", $true, $false, $false, $false, $false,
    $true, 1, $false, "               Payment ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "(£20 or £40 per survey: <br/>   see note below for details)"
#    -> "(see note below for details)      "
# ---------------------------------------------------------------------------
$searchPrefix = "    " + "   " + "(" + $pound + "20 or " + $pound + "40" + " per survey" + ": " + "^l" + "     "
$d.Content.Find.Execute($searchPrefix, $false, $false, $false, $false, $false,
    $true, 1, $false, "(", 2) | Out-Null

$d.Content.Find.Execute(
    "see note below for details)", $true, $false, $false, $false, $false,
    $true, 1, $false, "see note below for details)      ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "[Insert name of householder and address]"
#    -> "[Insert name of householder/ property owner & address]"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "[Insert name of householder and address]", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "[Insert name of householder/ property owner & address]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Rewrite the bold "Note: ..." paragraph.
# ---------------------------------------------------------------------------
$oldNote = "Note: You can claim payment for any CHEESE surveys that the householder has paid for that you have led (" + `
    $pound + "20 before CPD certificate; " + $pound + "40 after). For any CHEESE survey not paid for by the householder that you have led you can claim " + `
    $pound + "40 after CPD certification, but no payment can be claimed prior to CPD certification for these surveys due to the constraints of the project" + `
    $rsquo + "s finances"

$newNote = "Note: After your training surveys with Brian Harper and being signed off by Brian to lead surveys, you can claim payment of " + `
    $pound + "40 for any CHEESE surveys that you have led. For any CHEESE surveys that you assisted with, but did not lead, no payment can be claimed"

$d.Content.Find.Execute($oldNote, $true, $false, $false, $false, $false,
    $true, 1, $false, $newNote, 2) | Out-Null

# Underline just the "led" that now precedes ". For any CHEESE surveys..."
$rngLed = $d.Content
$rngLed.Find.Execute("led. For any CHEESE surveys that you assisted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ledOnly = $d.Range($rngLed.Start, $rngLed.Start + 3)
$ledOnly.Font.Underline = 1

# ---------------------------------------------------------------------------
# 6) Append the new sentence about expenses (with hyperlink) right after the
#    existing, unchanged trailing "." run of the Note paragraph.
# ---------------------------------------------------------------------------
$rngDot = $d.Content
$rngDot.Find.Execute("no payment can be claimed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $rngDot.End

$tail = $d.Range($insertAt, $insertAt)
$tail.InsertAfter(" ")
$boldSpace = $d.Range($insertAt, $insertAt + 1)
$boldSpace.Font.Bold = 1
$insertAt = $insertAt + 1

$sentence = "However, valid survey expenses " + $ndash + " such as car mileage " + $ndash + " can be claimed using the separate expenses form whether you led or assisted with a survey (expenses form available here: "
$tail2 = $d.Range($insertAt, $insertAt)
$tail2.InsertAfter($sentence)
$plainSentenceRange = $d.Range($insertAt, $insertAt + $sentence.Length)
$plainSentenceRange.Font.Bold = 0
$insertAt = $insertAt + $sentence.Length

$url = "https://cheeseproject.co.uk/energy-tracer-resources"
$tail3 = $d.Range($insertAt, $insertAt)
$tail3.InsertAfter($url)
$urlRange = $d.Range($insertAt, $insertAt + $url.Length)
$d.Hyperlinks.Add($urlRange, $url) | Out-Null
$insertAt = $insertAt + $url.Length

$closing = ")."
$tail4 = $d.Range($insertAt, $insertAt)
$tail4.InsertAfter($closing)
$closingRange = $d.Range($insertAt, $insertAt + $closing.Length)
$closingRange.Font.Bold = 0
